# LOQ4052.xlsx rebuild
#
# The published "syllabus" sheet had its content shifted by one row: a
# stray row (blank label in column A, with only the professor's name in
# columns B/C) sat at row 13. Removing that row realigns every label
# below it, and the custom row heights shift upward automatically with
# it. After the shift, a handful of B/C cells still hold stale text that
# needs to be corrected to match the refreshed build of the page.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray row 13 (blank A13; B13/C13 held the professor's name).
# This shifts rows 14:24 up to 13:23, carrying the row heights with them,
# and the sheet's used range shrinks from A1:C24 to A1:C23.
$ws.Rows(13).Delete()

# Fix up the handful of cells whose text fell out of sync with their
# labels after the shift.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# B15/C15 need the literal text "01/01/2011" (same as row 8), not an
# autoconverted date serial, so copy the already-text cell B8 and paste
# only its value into place.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

$ws.Range("B20").Value = "Provas e trabalhos"
$ws.Range("C20").Value = "Provas e trabalhos"

$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
